$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor has changed / signed off the timesheet: record her full name
# next to "Supervisor Name:" and her initials + sign-off date in the
# signature block (mirrors the existing employee sign-off row).
$ws.Range("G6").Value = "Prakruti Sinha"
$ws.Range("A27").Value = "P.S"

# Copy the date formatting from the employee's sign-off date cell so the
# new supervisor sign-off date cell (previously blank/unformatted) renders
# as a date like its sibling.
$ws.Range("D25").Copy()
$ws.Range("D27").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D27").Value = 41698        # 28/02/2014

# Reflect the new focus cell, as the author last had it selected.
$ws.Range("D27:E27").Select() | Out-Null
